$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows before the existing row 1246 (Excel shifts
# rows 1246:1285 down to 1248:1287, preserving their formatting/values).
$ws.Rows("1246:1247").Insert()

# Row 1246 - new "Primera" quality record for Brócoli, week of 2023-05-29
$ws.Cells.Item(1246, 1).Value = 6
$ws.Cells.Item(1246, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1246, 3).Value = "Metropolitana"
$ws.Cells.Item(1246, 4).Value = 45075
$ws.Cells.Item(1246, 5).Value = 13
$ws.Cells.Item(1246, 6).Value = 100112023
$ws.Cells.Item(1246, 7).Value = "Brócoli"
$ws.Cells.Item(1246, 8).Value = "Sin especificar"
$ws.Cells.Item(1246, 9).Value = "Primera"
$ws.Cells.Item(1246, 10).Value = 11400
$ws.Cells.Item(1246, 11).Value = 700
$ws.Cells.Item(1246, 12).Value = 800
$ws.Cells.Item(1246, 13).Value = 739
$ws.Cells.Item(1246, 14).Value = "$/unidad"
$ws.Cells.Item(1246, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1246, 16).Value = 739
$ws.Cells.Item(1246, 17).Value = 1
$ws.Cells.Item(1246, 18).Value = "Hortaliza"

# Row 1247 - new "Segunda" quality record for Brócoli, same week
$ws.Cells.Item(1247, 1).Value = 6
$ws.Cells.Item(1247, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1247, 3).Value = "Metropolitana"
$ws.Cells.Item(1247, 4).Value = 45075
$ws.Cells.Item(1247, 5).Value = 13
$ws.Cells.Item(1247, 6).Value = 100112023
$ws.Cells.Item(1247, 7).Value = "Brócoli"
$ws.Cells.Item(1247, 8).Value = "Sin especificar"
$ws.Cells.Item(1247, 9).Value = "Segunda"
$ws.Cells.Item(1247, 10).Value = 3800
$ws.Cells.Item(1247, 11).Value = 600
$ws.Cells.Item(1247, 12).Value = 600
$ws.Cells.Item(1247, 13).Value = 600
$ws.Cells.Item(1247, 14).Value = "$/unidad"
$ws.Cells.Item(1247, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1247, 16).Value = 600
$ws.Cells.Item(1247, 17).Value = 1
$ws.Cells.Item(1247, 18).Value = "Hortaliza"
